# Update entry to 2-bar confirmation and next-bar open
#
# Rewrites the trade log body (rows 2-35 of columns A-G) with the recomputed
# trade data that results from switching the strategy's entry rule to a
# 2-bar confirmation signal executed on the next bar's open. Every existing
# trade row is re-priced and nine new trades (rows 27-35) are appended at
# the end of the log, extending the sheet's used range from A1:G26 to
# A1:G35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Date(serial), Type, Entry, Exit, Close, PnL, Result
$data = @(
    @(45952, "LONG", 6678.3701171875, 6699.77978515625, 6699.77978515625, 21.40966796875, "Profit"),
    @(45953, "SHORT", 6747.14990234375, 6738.169921875, 6738.169921875, 8.97998046875, "Profit"),
    @(45954, "SHORT", 6804.06005859375, 6791.14013671875, 6791.14013671875, 12.919921875, "Profit"),
    @(45959, "LONG", 6899.18017578125, 6892.169921875, 6892.169921875, -7.01025390625, "Loss"),
    @(45961, "LONG", 6829.7001953125, 6840.22998046875, 6840.22998046875, 10.52978515625, "Profit"),
    @(45965, "LONG", 6781.39990234375, 6772.39990234375, 6772.39990234375, -9, "Loss"),
    @(45971, "SHORT", 6837.009765625, 6831.81982421875, 6831.81982421875, 5.18994140625, "Profit"),
    @(45972, "SHORT", 6850.68017578125, 6848.08984375, 6848.08984375, 2.59033203125, "Profit"),
    @(45974, "LONG", 6736.18994140625, 6738.52001953125, 6738.52001953125, 2.330078125, "Profit"),
    @(45978, "LONG", 6656.6201171875, 6672.91015625, 6672.91015625, 16.2900390625, "Profit"),
    @(45979, "SHORT", 6655.43994140625, 6616.93994140625, 6616.93994140625, 38.5, "Profit"),
    @(45980, "LONG", 6638.8701171875, 6642.93017578125, 6642.93017578125, 4.06005859375, "Profit"),
    @(45982, "SHORT", 6632.2001953125, 6603.7998046875, 6603.7998046875, 28.400390625, "Profit"),
    @(45985, "SHORT", 6700.14013671875, 6704.669921875, 6704.669921875, -4.52978515625, "Loss"),
    @(45986, "SHORT", 6747.3798828125, 6766.10009765625, 6766.10009765625, -18.72021484375, "Loss"),
    @(45992, "SHORT", 6819.5498046875, 6814.64990234375, 6814.64990234375, 4.89990234375, "Profit"),
    @(45994, "SHORT", 6858.22998046875, 6850.509765625, 6850.509765625, 7.72021484375, "Profit"),
    @(45999, "LONG", 6837.509765625, 6847.83984375, 6847.83984375, 10.330078125, "Profit"),
    @(46001, "SHORT", 6838.14990234375, 6886.9501953125, 6886.9501953125, -48.80029296875, "Loss"),
    @(46002, "SHORT", 6895.1298828125, 6902.10009765625, 6902.10009765625, -6.97021484375, "Loss"),
    @(46006, "LONG", 6817.14990234375, 6816.33984375, 6816.33984375, -0.81005859375, "Loss"),
    @(46007, "LONG", 6787.52978515625, 6798.509765625, 6798.509765625, 10.97998046875, "Profit"),
    @(46008, "LONG", 6736.6201171875, 6722.22021484375, 6722.22021484375, -14.39990234375, "Loss"),
    @(46009, "LONG", 6784.990234375, 6774.35009765625, 6774.35009765625, -10.64013671875, "Loss"),
    @(46010, "SHORT", 6831.52978515625, 6834.7001953125, 6834.7001953125, -3.17041015625, "Loss"),
    @(46014, "SHORT", 6905.7900390625, 6909.02001953125, 6909.02001953125, -3.22998046875, "Loss"),
    @(46017, "LONG", 6930.919921875, 6930.52001953125, 6930.52001953125, -0.39990234375, "Loss"),
    @(46024, "LONG", 6850.08984375, 6858.5400390625, 6858.5400390625, 8.4501953125, "Profit"),
    @(46028, "SHORT", 6944.740234375, 6946.330078125, 6946.330078125, -1.58984375, "Loss"),
    @(46031, "SHORT", 6973.68017578125, 6965.990234375, 6965.990234375, 7.68994140625, "Profit"),
    @(46034, "SHORT", 6978.7998046875, 6976.72021484375, 6976.72021484375, 2.07958984375, "Profit"),
    @(46035, "LONG", 6955.830078125, 6962.97998046875, 6962.97998046875, 7.14990234375, "Profit"),
    @(46036, "LONG", 6893.08984375, 6928.41015625, 6928.41015625, 35.3203125, "Profit"),
    @(46037, "SHORT", 6967.9599609375, 6944.97021484375, 6944.97021484375, 22.98974609375, "Profit")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

Write-Output "Updated $($data.Count) trade rows (A2:G$($startRow + $data.Count - 1))"
